$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.044.61'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.760.59'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.82'
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.18'
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -3.29%  '
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.84'
$ws.Range("E11").Value = '  -14.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  -3.66%  '
$ws.Range("D13").Value = '3.249.97'
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = '63.650.41'
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("E16").Value = '  -5.94%  '
$ws.Range("D17").Value = '2.767.92'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.10'
$ws.Range("E18").Value = '  -2.96%  '
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '359.42'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.63'
$ws.Range("E21").Value = '  -6.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.528'
$ws.Range("E23").Value = '  -8.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.02'
$ws.Range("E24").Value = '  -3.85%  '
$ws.Range("E25").Value = '  -3.98%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.51'
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").Value = '0.0₃0903'
$ws.Range("E28").Value = '  -7.16%  '
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("E30").Value = '  -4.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.28'
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '170.53'
$ws.Range("E32").Value = '  -1.40%  '
$ws.Range("E33").Value = '  -3.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.90'
$ws.Range("E34").Value = '  -4.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("E38").Value = '  -2.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '347.42'
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.17'
$ws.Range("E41").Value = '  -3.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.10'
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("E43").Value = '  -4.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.84'
$ws.Range("E44").Value = '  -4.12%  '
$ws.Range("E45").Value = '  -4.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '137.15'
$ws.Range("E46").Value = '  -1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.629'
$ws.Range("E47").Value = '  -3.77%  '
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.03'
$ws.Range("E51").Value = '  -0.03%  '
